$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Clear the material number / quantity cells that have now been filled via the
# new "fill material list" form helper, leaving the paired cell in each row intact.
$ws.Range("A12").ClearContents()
$ws.Range("B15").ClearContents()

# Move the active selection to B12, as left by the form-fill routine.
$ws.Range("B12").Select()
